$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (leading apostrophe) so numeric-looking strings
# such as "214.90" or "1.650.68" stay text cells instead of being
# auto-converted to numbers by Excel's smart-entry parsing.

$ws.Range('D2').Value = "'26.995.84"
$ws.Range('E2').Value = "'  +2.64%  "
$ws.Range('D3').Value = "'1.650.68"
$ws.Range('E3').Value = "'  +3.41%  "
$ws.Range('E4').Value = "'  +0.10%  "
$ws.Range('D5').Value = "'214.90"
$ws.Range('E5').Value = "'  +1.52%  "
$ws.Range('D6').Value = "'0.509"
$ws.Range('E6').Value = "'  +1.66%  "
$ws.Range('E7').Value = "'  +0.06%  "
$ws.Range('E8').Value = "'  +1.85%  "
$ws.Range('E9').Value = "'  +1.61%  "
$ws.Range('D10').Value = "'19.74"
$ws.Range('E10').Value = "'  +3.77%  "
$ws.Range('D11').Value = "'0.0867"
$ws.Range('E11').Value = "'  +1.46%  "
$ws.Range('D12').Value = "'1.884.93"
$ws.Range('E12').Value = "'  +3.51%  "
$ws.Range('D13').Value = "'1.653.21"
$ws.Range('E13').Value = "'  +3.64%  "
$ws.Range('D14').Value = "'4.08"
$ws.Range('E14').Value = "'  +2.40%  "
$ws.Range('D15').Value = "'0.518"
$ws.Range('E15').Value = "'  +2.98%  "
$ws.Range('D16').Value = "'65.25"
$ws.Range('E16').Value = "'  +2.81%  "
$ws.Range('D17').Value = "'239.40"
$ws.Range('E17').Value = "'  +4.41%  "
$ws.Range('D18').Value = "'26.977.47"
$ws.Range('E18').Value = "'  +2.58%  "
$ws.Range('D19').Value = "'7.86"
$ws.Range('E19').Value = "'  +2.68%  "
$ws.Range('E20').Value = "'  +1.17%  "
$ws.Range('E21').Value = "'  +0.10%  "
$ws.Range('E22').Value = "'  +4.46%  "
$ws.Range('D23').Value = "'2.25"
$ws.Range('E23').Value = "'  +4.04%  "
$ws.Range('D24').Value = "'9.24"
$ws.Range('E24').Value = "'  +3.54%  "
$ws.Range('D25').Value = "'145.60"
$ws.Range('E25').Value = "'  -0.52%  "
$ws.Range('E26').Value = "'  +0.00%  "
$ws.Range('D27').Value = "'7.10"
$ws.Range('E27').Value = "'  +1.99%  "
$ws.Range('E28').Value = "'  +2.12%  "
$ws.Range('D29').Value = "'15.81"
$ws.Range('E29').Value = "'  +2.64%  "
$ws.Range('D30').Value = "'0.0497"
$ws.Range('E30').Value = "'  +0.64%  "
$ws.Range('E31').Value = "'  +1.86%  "
$ws.Range('E32').Value = "'  +3.48%  "
$ws.Range('D33').Value = "'1.511.38"
$ws.Range('E33').Value = "'  +1.56%  "
$ws.Range('D34').Value = "'3.07"
$ws.Range('E34').Value = "'  +5.18%  "
$ws.Range('D35').Value = "'1.59"
$ws.Range('E35').Value = "'  +8.81%  "
$ws.Range('D36').Value = "'2.42"
$ws.Range('E36').Value = "'  -0.07%  "
$ws.Range('E37').Value = "'  +1.73%  "
$ws.Range('B38').Value = "'VeChain"
$ws.Range('C38').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('D38').Value = "'0.0169"
$ws.Range('E38').Value = "'  +2.76%  "
$ws.Range('B39').Value = "'ARBITRUM"
$ws.Range('C39').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D39').Value = "'0.884"
$ws.Range('E39').Value = "'  +8.38%  "
$ws.Range('D40').Value = "'5.95"
$ws.Range('E40').Value = "'  +2.72%  "
$ws.Range('E41').Value = "'  +0.03%  "
$ws.Range('B42').Value = "'MXToken"
$ws.Range('C42').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D42').Value = "'2.25"
$ws.Range('E42').Value = "'  +3.96%  "
$ws.Range('B43').Value = "'Aave"
$ws.Range('C43').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D43').Value = "'66.03"
$ws.Range('E43').Value = "'  +9.55%  "
$ws.Range('D44').Value = "'1.791.17"
$ws.Range('E44').Value = "'  +3.32%  "
$ws.Range('E45').Value = "'  +2.27%  "
$ws.Range('D46').Value = "'0.914"
$ws.Range('E46').Value = "'  -1.55%  "
$ws.Range('D47').Value = "'89.44"
$ws.Range('E47').Value = "'  +1.33%  "
$ws.Range('D48').Value = "'0.0₆0105"
$ws.Range('E48').Value = "'  +0.03%  "
$ws.Range('E49').Value = "'  +2.44%  "
$ws.Range('D50').Value = "'0.0507"
$ws.Range('E50').Value = "'  +1.30%  "
$ws.Range('D51').Value = "'0.0977"
$ws.Range('E51').Value = "'  +2.32%  "
